$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset push dropped two IPO listings (코칩 and SK증권제12호스팩) that
# previously occupied rows 20 and 21. Deleting the rows shifts the last
# row (신한글로벌액티브리츠, formerly row 22) up to row 20.
# Delete from the bottom up so row indices stay valid.
$ws.Range("A21:Y21").EntireRow.Delete()
$ws.Range("A20:Y20").EntireRow.Delete()
